$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 150.0354306666667
$ws.Range("H2").Value = 450.106292
$ws.Range("I2").Value = 0.4152507364956075
$ws.Range("J2").Value = 0.4152507364956075
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.607073666666667
$ws.Range("N2").Value = 25.821221
$ws.Range("O2").Value = 0.09431423806264486
$ws.Range("P2").Value = 0.09431423806264484
$ws.Range("Q2").Value = 1291.366004358059
$ws.Range("R2").Value = 11622.29403922253
$ws.Range("S2").Value = 0.03916405681753533
$ws.Range("T2").Value = 0.03916405681753533

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 150.0354306666667
$ws.Range("H3").Value = 450.106292
$ws.Range("I3").Value = 0.4152507364956075
$ws.Range("J3").Value = 0.4152507364956075
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 71.03134166666666
$ws.Range("N3").Value = 213.094025
$ws.Range("O3").Value = 0.7783443162342011
$ws.Range("P3").Value = 0.778344316234201
$ws.Range("Q3").Value = 10657.21793778948
$ws.Range("R3").Value = 95914.9614401053
$ws.Range("S3").Value = 0.3232080505634221
$ws.Range("T3").Value = 0.3232080505634219

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 150.0354306666667
$ws.Range("H4").Value = 450.106292
$ws.Range("I4").Value = 0.4152507364956075
$ws.Range("J4").Value = 0.4152507364956075
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 11.62112133333333
$ws.Range("N4").Value = 34.863364
$ws.Range("O4").Value = 0.1273414457031541
$ws.Range("P4").Value = 0.1273414457031541
$ws.Range("Q4").Value = 1743.579944076254
$ws.Range("R4").Value = 15692.21949668629
$ws.Range("S4").Value = 0.05287862911465015
$ws.Range("T4").Value = 0.05287862911465014

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 68.382243
$ws.Range("H5").Value = 205.146729
$ws.Range("I5").Value = 0.1892604742946246
$ws.Range("J5").Value = 0.1892604742946246
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.607073666666667
$ws.Range("N5").Value = 25.821221
$ws.Range("O5").Value = 0.09431423806264486
$ws.Range("P5").Value = 0.09431423806264484
$ws.Range("Q5").Value = 588.571002992901
$ws.Range("R5").Value = 5297.139026936108
$ws.Range("S5").Value = 0.01784995742847231
$ws.Range("T5").Value = 0.0178499574284723

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 68.382243
$ws.Range("H6").Value = 205.146729
$ws.Range("I6").Value = 0.1892604742946246
$ws.Range("J6").Value = 0.1892604742946246
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 71.03134166666666
$ws.Range("N6").Value = 213.094025
$ws.Range("O6").Value = 0.7783443162342011
$ws.Range("P6").Value = 0.778344316234201
$ws.Range("Q6").Value = 4857.282466466025
$ws.Range("R6").Value = 43715.54219819422
$ws.Range("S6").Value = 0.1473098144550102
$ws.Range("T6").Value = 0.1473098144550102

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 68.382243
$ws.Range("H7").Value = 205.146729
$ws.Range("I7").Value = 0.1892604742946246
$ws.Range("J7").Value = 0.1892604742946246
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 11.62112133333333
$ws.Range("N7").Value = 34.863364
$ws.Range("O7").Value = 0.1273414457031541
$ws.Range("P7").Value = 0.1273414457031541
$ws.Range("Q7").Value = 794.678342948484
$ws.Range("R7").Value = 7152.105086536355
$ws.Range("S7").Value = 0.02410070241114213
$ws.Range("T7").Value = 0.02410070241114213

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 104.737245
$ws.Range("H8").Value = 314.211735
$ws.Range("I8").Value = 0.2898796499701289
$ws.Range("J8").Value = 0.2898796499701289
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.607073666666667
$ws.Range("N8").Value = 25.821221
$ws.Range("O8").Value = 0.09431423806264486
$ws.Range("P8").Value = 0.09431423806264484
$ws.Range("Q8").Value = 901.4811833587149
$ws.Range("R8").Value = 8113.330650228434
$ws.Range("S8").Value = 0.0273397783167989
$ws.Range("T8").Value = 0.02733977831679889

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 104.737245
$ws.Range("H9").Value = 314.211735
$ws.Range("I9").Value = 0.2898796499701289
$ws.Range("J9").Value = 0.2898796499701289
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 71.03134166666666
$ws.Range("N9").Value = 213.094025
$ws.Range("O9").Value = 0.7783443162342011
$ws.Range("P9").Value = 0.778344316234201
$ws.Range("Q9").Value = 7439.627034820373
$ws.Range("R9").Value = 66956.64331338336
$ws.Range("S9").Value = 0.2256261779462096
$ws.Range("T9").Value = 0.2256261779462095

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 104.737245
$ws.Range("H10").Value = 314.211735
$ws.Range("I10").Value = 0.2898796499701289
$ws.Range("J10").Value = 0.2898796499701289
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 11.62112133333333
$ws.Range("N10").Value = 34.863364
$ws.Range("O10").Value = 0.1273414457031541
$ws.Range("P10").Value = 0.1273414457031541
$ws.Range("Q10").Value = 1217.16423226406
$ws.Range("R10").Value = 10954.47809037654
$ws.Range("S10").Value = 0.03691369370712048
$ws.Range("T10").Value = 0.03691369370712048

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 38.15794
$ws.Range("H11").Value = 114.47382
$ws.Range("I11").Value = 0.105609139239639
$ws.Range("J11").Value = 0.105609139239639
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.607073666666667
$ws.Range("N11").Value = 25.821221
$ws.Range("O11").Value = 0.09431423806264486
$ws.Range("P11").Value = 0.09431423806264484
$ws.Range("Q11").Value = 328.4282005482467
$ws.Range("R11").Value = 2955.85380493422
$ws.Range("S11").Value = 0.009960445499838322
$ws.Range("T11").Value = 0.009960445499838318

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 38.15794
$ws.Range("H12").Value = 114.47382
$ws.Range("I12").Value = 0.105609139239639
$ws.Range("J12").Value = 0.105609139239639
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 71.03134166666666
$ws.Range("N12").Value = 213.094025
$ws.Range("O12").Value = 0.7783443162342011
$ws.Range("P12").Value = 0.778344316234201
$ws.Range("Q12").Value = 2710.409673436167
$ws.Range("R12").Value = 24393.6870609255
$ws.Range("S12").Value = 0.08220027326955937
$ws.Range("T12").Value = 0.08220027326955932

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 38.15794
$ws.Range("H13").Value = 114.47382
$ws.Range("I13").Value = 0.105609139239639
$ws.Range("J13").Value = 0.105609139239639
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 11.62112133333333
$ws.Range("N13").Value = 34.863364
$ws.Range("O13").Value = 0.1273414457031541
$ws.Range("P13").Value = 0.1273414457031541
$ws.Range("Q13").Value = 443.4380505700533
$ws.Range("R13").Value = 3990.94245513048
$ws.Range("S13").Value = 0.01344842047024133
$ws.Range("T13").Value = 0.01344842047024133
